$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1247706666666667
$ws.Range("H2").Value = 0.374312
$ws.Range("I2").Value = 0.08987976235813114
$ws.Range("J2").Value = 0.08987976235813115
$ws.Range("M2").Value = 4.901461666666667
$ws.Range("N2").Value = 14.704385
$ws.Range("O2").Value = 0.2124427850531459
$ws.Range("P2").Value = 0.2124427850531459
$ws.Range("Q2").Value = 0.6115586397911111
$ws.Range("R2").Value = 5.504027758119999
$ws.Range("S2").Value = 0.01909430703527628
$ws.Range("T2").Value = 0.01909430703527629
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1247706666666667
$ws.Range("H3").Value = 0.374312
$ws.Range("I3").Value = 0.08987976235813114
$ws.Range("J3").Value = 0.08987976235813115
$ws.Range("O3").Value = 0.1372144215401173
$ws.Range("P3").Value = 0.1372144215401173
$ws.Range("Q3").Value = 0.3949988933528889
$ws.Range("R3").Value = 3.554990040176
$ws.Range("S3").Value = 0.01233279960013417
$ws.Range("T3").Value = 0.01233279960013417
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1247706666666667
$ws.Range("H4").Value = 0.374312
$ws.Range("I4").Value = 0.08987976235813114
$ws.Range("J4").Value = 0.08987976235813115
$ws.Range("M4").Value = 1.206743666666667
$ws.Range("N4").Value = 3.620231
$ws.Range("O4").Value = 0.05230357857032003
$ws.Range("P4").Value = 0.05230357857032004
$ws.Range("Q4").Value = 0.1505662117857778
$ws.Range("R4").Value = 1.355095906072
$ws.Range("S4").Value = 0.004701033212380205
$ws.Range("T4").Value = 0.004701033212380206
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1247706666666667
$ws.Range("H5").Value = 0.374312
$ws.Range("I5").Value = 0.08987976235813114
$ws.Range("J5").Value = 0.08987976235813115
$ws.Range("M5").Value = 13.79790933333333
$ws.Range("N5").Value = 41.393728
$ws.Range("O5").Value = 0.5980392148364168
$ws.Range("P5").Value = 0.5980392148364169
$ws.Range("Q5").Value = 1.721574346126222
$ws.Range("R5").Value = 15.494169115136
$ws.Range("S5").Value = 0.05375162251034047
$ws.Range("T5").Value = 0.05375162251034049
$ws.Range("I6").Value = 0.1613845890083672
$ws.Range("J6").Value = 0.1613845890083672
$ws.Range("M6").Value = 4.901461666666667
$ws.Range("N6").Value = 14.704385
$ws.Range("O6").Value = 0.2124427850531459
$ws.Range("P6").Value = 0.2124427850531459
$ws.Range("Q6").Value = 1.098090795388889
$ws.Range("R6").Value = 9.8828171585
$ws.Range("S6").Value = 0.03428499155359484
$ws.Range("T6").Value = 0.03428499155359485
$ws.Range("I7").Value = 0.1613845890083672
$ws.Range("J7").Value = 0.1613845890083672
$ws.Range("O7").Value = 0.1372144215401173
$ws.Range("P7").Value = 0.1372144215401173
$ws.Range("S7").Value = 0.02214429302627268
$ws.Range("T7").Value = 0.02214429302627268
$ws.Range("I8").Value = 0.1613845890083672
$ws.Range("J8").Value = 0.1613845890083672
$ws.Range("M8").Value = 1.206743666666667
$ws.Range("N8").Value = 3.620231
$ws.Range("O8").Value = 0.05230357857032003
$ws.Range("P8").Value = 0.05230357857032004
$ws.Range("Q8").Value = 0.2703508061222222
$ws.Range("R8").Value = 2.4331572551
$ws.Range("S8").Value = 0.008440991531237941
$ws.Range("T8").Value = 0.008440991531237943
$ws.Range("I9").Value = 0.1613845890083672
$ws.Range("J9").Value = 0.1613845890083672
$ws.Range("M9").Value = 13.79790933333333
$ws.Range("N9").Value = 41.393728
$ws.Range("O9").Value = 0.5980392148364168
$ws.Range("P9").Value = 0.5980392148364169
$ws.Range("Q9").Value = 3.091191620977777
$ws.Range("R9").Value = 27.8207245888
$ws.Range("S9").Value = 0.09651431289726174
$ws.Range("T9").Value = 0.09651431289726177
$ws.Range("G10").Value = 1.039391333333333
$ws.Range("H10").Value = 3.118174
$ws.Range("I10").Value = 0.7487356486335016
$ws.Range("J10").Value = 0.7487356486335016
$ws.Range("M10").Value = 4.901461666666667
$ws.Range("N10").Value = 14.704385
$ws.Range("O10").Value = 0.2124427850531459
$ws.Range("P10").Value = 0.2124427850531459
$ws.Range("Q10").Value = 5.094536776998889
$ws.Range("R10").Value = 45.85083099299
$ws.Range("S10").Value = 0.1590634864642747
$ws.Range("T10").Value = 0.1590634864642748
$ws.Range("G11").Value = 1.039391333333333
$ws.Range("H11").Value = 3.118174
$ws.Range("I11").Value = 0.7487356486335016
$ws.Range("J11").Value = 0.7487356486335016
$ws.Range("O11").Value = 0.1372144215401173
$ws.Range("P11").Value = 0.1372144215401173
$ws.Range("Q11").Value = 3.290504390139112
$ws.Range("R11").Value = 29.614539511252
$ws.Range("S11").Value = 0.1027373289137104
$ws.Range("T11").Value = 0.1027373289137104
$ws.Range("G12").Value = 1.039391333333333
$ws.Range("H12").Value = 3.118174
$ws.Range("I12").Value = 0.7487356486335016
$ws.Range("J12").Value = 0.7487356486335016
$ws.Range("M12").Value = 1.206743666666667
$ws.Range("N12").Value = 3.620231
$ws.Range("O12").Value = 0.05230357857032003
$ws.Range("P12").Value = 0.05230357857032004
$ws.Range("Q12").Value = 1.254278908688222
$ws.Range("R12").Value = 11.288510178194
$ws.Range("S12").Value = 0.03916155382670188
$ws.Range("T12").Value = 0.03916155382670189
$ws.Range("G13").Value = 1.039391333333333
$ws.Range("H13").Value = 3.118174
$ws.Range("I13").Value = 0.7487356486335016
$ws.Range("J13").Value = 0.7487356486335016
$ws.Range("M13").Value = 13.79790933333333
$ws.Range("N13").Value = 41.393728
$ws.Range("O13").Value = 0.5980392148364168
$ws.Range("P13").Value = 0.5980392148364169
$ws.Range("Q13").Value = 14.34142737918578
$ws.Range("R13").Value = 129.072846412672
$ws.Range("S13").Value = 0.4477732794288145
$ws.Range("T13").Value = 0.4477732794288146
